# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets
# to reflect newly generated output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4746
$ws1.Range("F4").Value = 3208
$ws1.Range("F5").Value = 670
$ws1.Range("F11").Value = 1411
$ws1.Range("F13").Value = 1677
$ws1.Range("F15").Value = 135
$ws1.Range("F21").Value = 65
$ws1.Range("F22").Value = 127
$ws1.Range("F23").Value = 15
$ws1.Range("F25").Value = 53
$ws1.Range("F27").Value = 4232
$ws1.Range("F28").Value = 20
$ws1.Range("F29").Value = 789
$ws1.Range("F31").Value = 2051
$ws1.Range("F33").Value = 1963

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4746
$ws4.Range("F4").Value = 3208
$ws4.Range("F5").Value = 670
$ws4.Range("F12").Value = 1411
$ws4.Range("F14").Value = 1677
$ws4.Range("F16").Value = 135
$ws4.Range("F22").Value = 65
$ws4.Range("F23").Value = 127
$ws4.Range("F24").Value = 15
$ws4.Range("F26").Value = 53
$ws4.Range("F28").Value = 4232
$ws4.Range("F30").Value = 20
$ws4.Range("F32").Value = 789
$ws4.Range("F34").Value = 2051
$ws4.Range("F36").Value = 1963
